$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-15 Tuesday" "2024-10-16 Wednesday"

Replace-Text "544×5=2720" "670×8=5360"
Replace-Text "719×5=3595" "206×5=1030"
Replace-Text "829×4=3316" "427×7=2989"
Replace-Text "830×8=6640" "847×2=1694"
Replace-Text "573×3=1719" "953×3=2859"
Replace-Text "920×3=2760" "519×5=2595"
Replace-Text "426×8=3408" "265×9=2385"
Replace-Text "803×6=4818" "843×9=7587"
Replace-Text "339×6=2034" "155×3=465"
Replace-Text "635×8=5080" "365×3=1095"
Replace-Text "626×7=4382" "657×6=3942"
Replace-Text "278×9=2502" "626×4=2504"
Replace-Text "864×2=1728" "931×6=5586"
Replace-Text "116×4=464" "970×6=5820"
Replace-Text "813×3=2439" "585×3=1755"
Replace-Text "814×2=1628" "113×9=1017"
Replace-Text "694×3=2082" "443×6=2658"
Replace-Text "142×6=852" "321×6=1926"
Replace-Text "798×2=1596" "187×3=561"
Replace-Text "250×2=500" "799×7=5593"
Replace-Text "974×3=2922" "198×8=1584"
Replace-Text "140×3=420" "314×8=2512"
Replace-Text "770×2=1540" "971×5=4855"
Replace-Text "139×6=834" "588×5=2940"
Replace-Text "897×4=3588" "866×6=5196"
